# Auto-generated Excel COM-interop script to apply Marilith_Profits market-price updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1048.45
$ws.Range("I32").Value = 816.36365
$ws.Range("K32").Value = 816.36365
$ws.Range("M32").Value = -490.36365
$ws.Range("H33").Value = 180.73334
$ws.Range("I33").Value = 189.63637
$ws.Range("J33").Value = 156.25
$ws.Range("K33").Value = 189.63637
$ws.Range("L33").Value = 156.25
$ws.Range("M33").Value = 39.36363
$ws.Range("N33").Value = -614.25
$ws.Range("H105").Value = 34583.332
$ws.Range("J105").Value = 34583.332
$ws.Range("L105").Value = 34583.332
$ws.Range("N105").Value = -41571.332
$ws.Range("H132").Value = 6333.6665
$ws.Range("J132").Value = 2500
$ws.Range("L132").Value = 7500
$ws.Range("N132").Value = -12560
$ws.Range("H137").Value = 3570.1428
$ws.Range("J137").Value = 4198.2
$ws.Range("L137").Value = 12594.6
$ws.Range("N137").Value = -17694.6
$ws.Range("H138").Value = 3571.4285
$ws.Range("J138").Value = 3571.4285
$ws.Range("L138").Value = 10714.2855
$ws.Range("N138").Value = -20994.2855
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2018.9
$ws.Range("I2").Value = 997.8
$ws.Range("K2").Value = 997.8
$ws.Range("M2").Value = -884.8
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H37").Value = 20000
$ws.Range("J37").Value = 20000
$ws.Range("L37").Value = 20000
$ws.Range("N37").Value = -20546
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H116").Value = 2018.9
$ws.Range("I116").Value = 997.8
$ws.Range("K116").Value = 997.8
$ws.Range("M116").Value = 1296.2
$ws.Range("H118").Value = 63666
$ws.Range("J118").Value = 63666
$ws.Range("L118").Value = 63666
$ws.Range("N118").Value = -66980
$ws.Range("H132").Value = 1769.72
$ws.Range("I132").Value = 1776.8334
$ws.Range("J132").Value = 1599
$ws.Range("K132").Value = 5330.5002
$ws.Range("L132").Value = 4797
$ws.Range("M132").Value = -2800.5002
$ws.Range("N132").Value = -9857
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2018.9
$ws.Range("I3").Value = 997.8
$ws.Range("K3").Value = 997.8
$ws.Range("M3").Value = -883.8
$ws.Range("H19").Value = 4750
$ws.Range("I19").Value = 3500
$ws.Range("J19").Value = 6000
$ws.Range("K19").Value = 3500
$ws.Range("L19").Value = 6000
$ws.Range("M19").Value = -3327
$ws.Range("N19").Value = -6346
$ws.Range("H94").Value = 2620.5715
$ws.Range("I94").Value = 2668.611
$ws.Range("J94").Value = 2332.3333
$ws.Range("K94").Value = 2668.611
$ws.Range("L94").Value = 2332.3333
$ws.Range("M94").Value = -2217.611
$ws.Range("N94").Value = -3234.3333
$ws.Range("H134").Value = 3598.6667
$ws.Range("I134").Value = 3650.8262
$ws.Range("K134").Value = 10952.4786
$ws.Range("M134").Value = -8417.4786
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 2002673.8
$ws.Range("I3").Value = 3335040
$ws.Range("J3").Value = 4124.5
$ws.Range("K3").Value = 3335040
$ws.Range("L3").Value = 4124.5
$ws.Range("M3").Value = -3334927
$ws.Range("N3").Value = -4350.5
$ws.Range("H22").Value = 995
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 995
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 995
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1695
$ws.Range("H33").Value = 1302
$ws.Range("I33").Value = 1302
$ws.Range("K33").Value = 1302
$ws.Range("M33").Value = -923
$ws.Range("H58").Value = 3584.5
$ws.Range("I58").Value = 4569
$ws.Range("J58").Value = 2600
$ws.Range("K58").Value = 4569
$ws.Range("L58").Value = 2600
$ws.Range("M58").Value = -4366
$ws.Range("N58").Value = -3006
$ws.Range("H134").Value = 1563.3182
$ws.Range("I134").Value = 1572.15
$ws.Range("J134").Value = 1475
$ws.Range("K134").Value = 4716.450000000001
$ws.Range("L134").Value = 4425
$ws.Range("M134").Value = -2181.450000000001
$ws.Range("N134").Value = -9495
$ws.Range("H136").Value = 3584.5
$ws.Range("I136").Value = 4569
$ws.Range("J136").Value = 2600
$ws.Range("K136").Value = 13707
$ws.Range("L136").Value = 7800
$ws.Range("M136").Value = -11157
$ws.Range("N136").Value = -12900
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 5000
$ws.Range("I57").Value = 5000
$ws.Range("K57").Value = 15000
$ws.Range("M57").Value = -14441
$ws.Range("H109").Value = 445
$ws.Range("I109").Value = 445
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 1335
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -295
$ws.Range("N109").ClearContents()
$ws.Range("H117").Value = 753.5
$ws.Range("I117").Value = 755.3333
$ws.Range("K117").Value = 2265.9999
$ws.Range("M117").Value = 1176.0001
$ws.Range("H128").Value = 339792.88
$ws.Range("I128").Value = 339792.88
$ws.Range("K128").Value = 1019378.64
$ws.Range("M128").Value = -1014398.64
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = 0
$ws.Range("H62").Value = 27777.777
$ws.Range("I62").Value = 26000
$ws.Range("K62").Value = 26000
$ws.Range("M62").Value = -25314
$ws.Range("H65").Value = 27777.777
$ws.Range("I65").Value = 26000
$ws.Range("K65").Value = 78000
$ws.Range("M65").Value = -74568
$ws.Range("H80").Value = 3747.5
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 4197
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 4197
$ws.Range("M80").Value = -502
$ws.Range("N80").Value = -6193
$ws.Range("H83").Value = 3747.5
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 4197
$ws.Range("K83").Value = 7500
$ws.Range("L83").Value = 20985
$ws.Range("M83").Value = -2508
$ws.Range("N83").Value = -30969
$ws.Range("H102").Value = 2549.25
$ws.Range("I102").Value = 2399
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2399
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -777
$ws.Range("N102").Value = -6244
$ws.Range("H122").Value = 3939.5557
$ws.Range("I122").Value = 2491.3333
$ws.Range("K122").Value = 7473.999899999999
$ws.Range("M122").Value = -5023.999899999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H46").Value = 3846.9412
$ws.Range("I46").Value = 2733
$ws.Range("K46").Value = 2733
$ws.Range("M46").Value = -2545
$ws.Range("H63").Value = 29998
$ws.Range("H66").Value = 29998
$ws.Range("H93").Value = 890.8
$ws.Range("I93").Value = 890.8
$ws.Range("K93").Value = 890.8
$ws.Range("M93").Value = 357.2
$ws.Range("H100").Value = 1320.6
$ws.Range("I100").Value = 1350.75
$ws.Range("K100").Value = 1350.75
$ws.Range("M100").Value = -809.75
$ws.Range("H114").Value = 33499.668
$ws.Range("J114").Value = 33499.668
$ws.Range("L114").Value = 33499.668
$ws.Range("N114").Value = -42177.668
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 10000
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H136").Value = 13995
$ws.Range("I136").Value = 13995
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 41985
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -39435
$ws.Range("N136").ClearContents()
